$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.827.35"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.385.09"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'557.15"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "'133.49"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.584"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").Value = "'5.61"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").Value = "'24.40"
$ws.Range("E13").Value = "  -4.40%  "
$ws.Range("D14").Value = "2.815.30"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "59.786.63"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "'0.0000136"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "2.391.86"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "'11.11"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").Value = "'320.57"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "'6.72"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'64.18"
$ws.Range("E23").Value = "  -3.62%  "
$ws.Range("D24").Value = "'0.172"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'8.43"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'1.80"
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").Value = "0.0₃0759"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("D30").Value = "'170.25"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "'6.05"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("E32").Value = "  +5.66%  "
$ws.Range("D33").Value = "'0.398"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").Value = "'18.15"
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'4.13"
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("D39").Value = "'1.58"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("D40").Value = "'318.48"
$ws.Range("D41").Value = "'38.65"
$ws.Range("D42").Value = "'146.39"
$ws.Range("E42").Value = "  +5.49%  "
$ws.Range("D43").Value = "'3.52"
$ws.Range("E43").Value = "  -4.15%  "
$ws.Range("D44").Value = "'0.0966"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "'19.69"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").Value = "'0.0510"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Value = "'0.571"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").Value = "'0.0217"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").Value = "'11.06"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("B50").Value = "ZEEBU"
$ws.Range("C50").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D50").Value = "'4.68"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").Value = "'0.949"
$ws.Range("E51").Value = "  +0.11%  "
